# Update init file cell references for extra scenarios (26,28,32)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Column D holds the "Lower Right Cell" reference of several data blocks.
# Extend the blocks from row 40 to row 43 (and the Variables List block from E258 to E261).
$ws.Range("D5").Value  = "A43"
$ws.Range("D6").Value  = "B43"
$ws.Range("D7").Value  = "C43"
$ws.Range("D8").Value  = "G43"
$ws.Range("D9").Value  = "H43"
$ws.Range("D10").Value = "I43"
$ws.Range("D11").Value = "J43"
$ws.Range("D15").Value = "E261"

# D22 keeps the same text (O475) but make sure it stays set explicitly.
$ws.Range("D22").Value = "O475"

# Update the active selection to reflect where the user ended up working.
$ws.Range("D15").Select()
